# "my new running extent"
# Rename the sheets, add the two new data points (lakshma / Infosys),
# resize the "delete" sheet's first column, and move the active
# selection/tab so sheet2 ("delete") ends up the active tab with A2
# selected while sheet1 ("edit") is left with D14 selected.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets (Sheet1 -> edit, Sheet2 -> delete; Sheet3 unchanged).
$ws1.Name = "edit"
$ws2.Name = "delete"

# New content.
$ws1.Range("A3").Value = "lakshma"
$ws2.Range("A1").Value = "campaign name"
$ws2.Range("A2").Value = "Infosys"

# Column width on the "delete" sheet.
$ws2.Columns.Item(1).ColumnWidth = 12.7109375

# Selections on each sheet.
$ws1.Range("D14").Select() | Out-Null
$ws2.Range("A2").Select() | Out-Null

# Make "delete" the active tab (also moves tabSelected from sheet1 to sheet2).
$ws2.Activate() | Out-Null

Write-Output "done"
